$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Remove the extra (empty) sheets Tabelle2 and Tabelle3 ---
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Tabelle2").Delete() | Out-Null
$wb.Worksheets.Item("Tabelle3").Delete() | Out-Null

# --- Update existing comment text ---
$ws.Range("D1").Comment.Text("Language comment") | Out-Null
$ws.Range("A2").Comment.Text("Concept comment") | Out-Null

# --- Update B2 value (was "I", becomes "Woman") ---
$ws.Range("B2").Value = "Woman"

# --- Add new row 3 (keep the shared-string interning order matching the source data) ---
$ws.Range("A3").Value = "TG100"
$ws.Range("D3").Value = "/am/ (description) {anysource}"
$ws.Range("B3").Value = "Person"

# --- Give A3 the same style as A2 ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Add a comment on the new row's A3 cell ---
$ws.Range("A3").AddComment("Concept comment") | Out-Null

# --- Restore the active-cell selection ---
$ws.Range("D18").Select() | Out-Null

Write-Host "Edit complete"
